$d = $word.ActiveDocument

# The original sentence (single run) needs to become three runs:
#   "On the first graph, ... This could be because "
#   "these channels posted neutral comments"
#   ", or because Vader Analysis ... polarity. "
#
# We replace the middle portion "users post neutral comments about these channels"
# with "these channels posted neutral comments", splitting the run into three pieces
# by using Find/Replace on the specific substring while preserving the surrounding text
# as separate runs (achieved by selecting and replacing just that middle substring,
# which naturally splits the run in Word's editing model).

$rng = $d.Content
$rng.Find.Execute("users post neutral comments about these channels", $true, $false, $false, $false, $false, $true, 1, $false, "these channels posted neutral comments", 2) | Out-Null
